$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Add the "Errors" sheet after Sheet1
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "Errors"

$newSheet.Range('A1').Value = 'COLUMN'
$newSheet.Range('B1').Value = 'ROW'
$newSheet.Range('C1').Value = 'Error'
$newSheet.Range('D1').Value = 'Expected'
$newSheet.Range('E1').Value = 'Recieved'
$newSheet.Range('F1').Value = 'Warn If Error Test'
$newSheet.Range('G1').Value = 'Empty'
$newSheet.Range('H1').Value = 'Spaces'
$newSheet.Range('I1').Value = 'None-Alphanumeric'
$newSheet.Range('J1').Value = 'contains Numbers'
$newSheet.Range('K1').Value = 'Contains Letters'
$newSheet.Range('A2').Value = 'Spaces'
$newSheet.Range('B2').Value = 2
$newSheet.Range('C2').Value = 'Cell Contains a Space'
$newSheet.Range('D2').Value = 'No Spaces'
$newSheet.Range('E2').Value = '     hbh    hgvjh'
$newSheet.Range('F2').Value = $true
$newSheet.Range('H2').Value = '     hbh    hgvjh'
$newSheet.Range('I2').Value = '/// &&'
$newSheet.Range('J2').Value = 5555
$newSheet.Range('K2').Value = 'sds DFSDF'
$newSheet.Range('A3').Value = 'None-Alphanumeric'
$newSheet.Range('B3').Value = 2
$newSheet.Range('C3').Value = 'Cell Contains Non-AlphaNumerics'
$newSheet.Range('D3').Value = 'A-Z / 0-9'
$newSheet.Range('E3').Value = '/// &&'
$newSheet.Range('F3').Value = $true
$newSheet.Range('H3').Value = '     hbh    hgvjh'
$newSheet.Range('I3').Value = '/// &&'
$newSheet.Range('J3').Value = 5555
$newSheet.Range('K3').Value = 'sds DFSDF'
$newSheet.Range('A4').Value = 'None-Alphanumeric'
$newSheet.Range('B4').Value = 2
$newSheet.Range('C4').Value = 'Cell Contains Non-AlphaNumerics'
$newSheet.Range('D4').Value = 'A-Z / 0-9'
$newSheet.Range('E4').Value = '/// &&'
$newSheet.Range('F4').Value = $true
$newSheet.Range('H4').Value = '     hbh    hgvjh'
$newSheet.Range('I4').Value = '/// &&'
$newSheet.Range('J4').Value = 5555
$newSheet.Range('K4').Value = 'sds DFSDF'
$newSheet.Range('A5').Value = 'None-Alphanumeric'
$newSheet.Range('B5').Value = 2
$newSheet.Range('C5').Value = 'Cell Contains Non-AlphaNumerics'
$newSheet.Range('D5').Value = 'A-Z / 0-9'
$newSheet.Range('E5').Value = '/// &&'
$newSheet.Range('F5').Value = $true
$newSheet.Range('H5').Value = '     hbh    hgvjh'
$newSheet.Range('I5').Value = '/// &&'
$newSheet.Range('J5').Value = 5555
$newSheet.Range('K5').Value = 'sds DFSDF'
$newSheet.Range('A6').Value = 'None-Alphanumeric'
$newSheet.Range('B6').Value = 2
$newSheet.Range('C6').Value = 'Cell Contains Non-AlphaNumerics'
$newSheet.Range('D6').Value = 'A-Z / 0-9'
$newSheet.Range('E6').Value = '/// &&'
$newSheet.Range('F6').Value = $true
$newSheet.Range('H6').Value = '     hbh    hgvjh'
$newSheet.Range('I6').Value = '/// &&'
$newSheet.Range('J6').Value = 5555
$newSheet.Range('K6').Value = 'sds DFSDF'
$newSheet.Range('A7').Value = 'None-Alphanumeric'
$newSheet.Range('B7').Value = 2
$newSheet.Range('C7').Value = 'Cell Contains Non-AlphaNumerics'
$newSheet.Range('D7').Value = 'A-Z / 0-9'
$newSheet.Range('E7').Value = '/// &&'
$newSheet.Range('F7').Value = $true
$newSheet.Range('H7').Value = '     hbh    hgvjh'
$newSheet.Range('I7').Value = '/// &&'
$newSheet.Range('J7').Value = 5555
$newSheet.Range('K7').Value = 'sds DFSDF'
$newSheet.Range('A8').Value = 'None-Alphanumeric'
$newSheet.Range('B8').Value = 2
$newSheet.Range('C8').Value = 'Cell Contains Non-AlphaNumerics'
$newSheet.Range('D8').Value = 'A-Z / 0-9'
$newSheet.Range('E8').Value = '/// &&'
$newSheet.Range('F8').Value = $true
$newSheet.Range('H8').Value = '     hbh    hgvjh'
$newSheet.Range('I8').Value = '/// &&'
$newSheet.Range('J8').Value = 5555
$newSheet.Range('K8').Value = 'sds DFSDF'
$newSheet.Range('A9').Value = 'Contains Letters'
$newSheet.Range('B9').Value = 2
$newSheet.Range('C9').Value = 'Cell Contains Letters'
$newSheet.Range('D9').Value = 'String'
$newSheet.Range('E9').Value = 'sds DFSDF'
$newSheet.Range('F9').Value = $true
$newSheet.Range('H9').Value = '     hbh    hgvjh'
$newSheet.Range('I9').Value = '/// &&'
$newSheet.Range('J9').Value = 5555
$newSheet.Range('K9').Value = 'sds DFSDF'
# Header row fill color
$newSheet.Range("A1:K1").Interior.Color = 9145088

Write-Host "Done building Errors sheet"
